$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the existing
# header style (bold, bordered, centered) from H1 so the new columns
# look the same as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-7
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 3

$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 5
